$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# The sheet currently ends with a data row (104) followed by a footnote
# row (105, "*since 4/8 the two hotlines were merged"). A new day's
# figures (2020-05-09) need to be appended as row 105, which pushes the
# footnote row down to 106. Inserting a row above the footnote row picks
# up the same number formats/styles used by the rest of the data rows.
$ws.Rows.Item(105).Insert()

$ws.Range("A105").Value = 43960
$ws.Range("B105").Value = 378
$ws.Range("C105").Value = 35385
$ws.Range("D105").Value = 98
$ws.Range("E105").Value = 7232

$ws.Range("E106").Select()

# Grow the sheet's print area to cover the newly added row. Drop the
# existing (workbook-level) Print_Area definition and recreate it scoped
# to this worksheet, as Excel stores `_xlnm.Print_Area`.
for ($i = $wb.Names.Count(); $i -ge 1; $i--) {
  $n = $wb.Names.Item($i)
  if ($n.Name() -like "*Print_Area*") {
    $n.Delete()
  }
}
$ws.Names.Add("_xlnm.Print_Area", "=相談件数!`$A`$1:`$E`$107")
